$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the Field values (column A) that used mixed-case / verbose naming
# to the new lower_snake_case convention. Only the rows whose field name
# actually changed are touched here (in the same order they were fixed).
$ws.Range("A44").Value = "nspike_ratio"
$ws.Range("A29").Value = "v_threshold"
$ws.Range("A31").Value = "first_ap_peak_amplitude"
$ws.Range("A32").Value = "first_ap_peak_time"
$ws.Range("A33").Value = "first_ap_trough_amplitude"
$ws.Range("A34").Value = "first_ap_trough_time"
$ws.Range("A39").Value = "max_isi_cv"
$ws.Range("A45").Value = "max_ahp_after_depol_injection"
$ws.Range("A24").Value = "resistance_rsquared"

# Column C (Description) is no longer used for rows 23-48; clear it out.
$ws.Range("C23:C48").ClearContents()

# Restore the selection / view state to match the saved workbook.
$ws.Range("A15").Select()
$ws.Range("C23:C48").Select()
